# H1AR20 BOM update: swap out the two LED parts.
#  D1: VLMS1300-GS08 (Red)         -> VLMY1300-GS08 (Yellow)
#  D2: Cree CLVBA-FKA-CC1F1L1BB7R3R3 (RGB) -> VLMO1300-GS08 (Soft Orange)
# Both new parts are Vishay, so the MFN column for D2 also changes to Vishay.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 10 (designator D1): VLMS1300-GS08 -> VLMY1300-GS08 ---------------
$ws.Range("B10").Value = "VLMY1300-GS08"
$ws.Range("D10").Value = "LED Uni-Color Yellow 588nm 2-Pin Chip 0603(1608Metric) T/R"
$ws.Range("E10").Value = "Vishay"
$ws.Range("F10").Value = "VLMY1300-GS08"
$ws.Range("H10").Value = "https://octopart.com/vlmy1300-gs08-vishay-21709204?r=sp"

# --- Row 11 (designator D2): Cree CLVBA-FKA-CC1F1L1BB7R3R3 -> VLMO1300-GS08 ---
$ws.Range("B11").Value = "VLMO1300-GS08 "
$ws.Range("D11").Value = "LED Uni-Color Soft Orange 611nm 2-Pin Chip 0603(1608Metric) T/R"
$ws.Range("E11").Value = "Vishay"
$ws.Range("F11").Value = "VLMO1300-GS08 "
$ws.Range("H11").Value = "https://octopart.com/vlmo1300-gs08-vishay-21709200?r=sp#"

# D2's OCTOPART_URL cell had no real hyperlink before (just hyperlink-style
# formatting); the new part's URL gets an actual clickable hyperlink.
$ws.Hyperlinks.Add($ws.Range("H11"), "https://octopart.com/vlmo1300-gs08-vishay-21709200?r=sp#")

# Leave selection where the author left it after editing the D2 description.
$ws.Range("D11").Select()
